$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -5.724999999999999
$ws.Range("C9").Value = -11.99470000000001
$ws.Range("D11").Value = -8.531600000000001
$ws.Range("C18").Value = -14.31479999999999
$ws.Range("C20").Value = -13.53649999999998
$ws.Range("E21").Value = 13.07259999999999
